# Add data for 2022-07-04:
# The "June 2022" running-total column (column B) now reflects data
# through June 26 instead of June 25, so the sheet/tab name, the column
# header text, and a handful of neighborhood/month cell counts change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name == workbook.xml sheet name).
$ws.Name = "Through 2022-06-26"

# Update the column B header text shared by the sheet.
$ws.Range("B1").Value = "June 2022 (through June 26)"

# Updated / added carjacking counts.
$ws.Range("B2").Value = 8
$ws.Range("AL2").Value = 4
$ws.Range("AF5").Value = 9
$ws.Range("N7").Value = 3
$ws.Range("H8").Value = 1
$ws.Range("B9").Value = 6
$ws.Range("H9").Value = 4
$ws.Range("Z9").Value = 3
$ws.Range("AF11").Value = 1
$ws.Range("N12").Value = 4
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 8
$ws.Range("Z14").Value = 4
$ws.Range("N19").Value = 3
$ws.Range("AL23").Value = 2
$ws.Range("N33").Value = 1
$ws.Range("N34").Value = 1
$ws.Range("AR35").Value = 1
$ws.Range("T36").Value = 1
$ws.Range("AF38").Value = 1
$ws.Range("H70").Value = 3
$ws.Range("B81").Value = 3
$ws.Range("B90").Value = 1
$ws.Range("AF92").Value = 1
$ws.Range("B94").Value = 3
